$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# A2 uses style 1 -> fontId 2. Let's touch Name (same value) and see if new font created.
$ws.Range("A2").Font.Name = "等线"
